# Generate Report for Handoff
#
# Updates the "latest handoff" timestamps for the files that were just
# (re-)handed off (i.e. rows whose status is "Handback transform failed" or
# "Ready for handoff"), across the Overview sheet and each language sheet.
#
# Rows left untouched:
#   - row 2 (644e43b3...) "Handed back: in sync with en-US"
#   - row 3 (8e59d3a0...) "Handed back: in sync with en-US"
#   - row 5 (dec38479...) "In Translation"

$wb = $excel.ActiveWorkbook

# Rows that received a new handoff timestamp.
$rows = @(4, 6, 7, 8, 9, 10)

# Overview sheet: column D = "Latest Handoff Date"
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 4).Value = "2016-35-18 03:35:20"
}

# zh-cn sheet: column E = "Latest Handoff Datetime"
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 5).Value = "2016-03-18 03:35:11"
}

# de-de sheet: column E = "Latest Handoff Datetime"
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 5).Value = "2016-03-18 03:35:20"
}
